$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# The "Bought (feemicon pill) X 10000" row (row 13) is removed; product stock is
# now decreased on purchase/sale instead of being tracked as one giant bulk buy.
$ws1.Rows(13).Delete()

# Newly recorded buy/sell transactions appended to the report.
$newRows = @(
  @(44864.9072121412,    "Sold (PD Name 4) X 10",            150.0),
  @(44864.90856619213,   "Sold (PD Name 1) X 20",             240.0),
  @(44864.90928422454,   "Sold (PD Name 2) X 12",             156.0),
  @(44864.91438642361,   "Sold (PD Name 4) X 1",               15.0),
  @(44864.91845265046,   "Sold (PD Name 1) X 12",             144.0),
  @(44864.918554965276,  "Sold (PD Name 2) X 1",               13.0),
  @(44864.918612615744,  "Sold (PD Name 7) X 12",             216.0),
  @(44864.91877981481,   "Sold (PD Name 1) X 12",             144.0),
  @(44864.92317155092,   "Bought (Third Party PD - 3) X 6",   -84.0),
  @(44864.92766209491,   "Bought (Third Party PD - 1) X 10",  -120.0),
  @(44864.92783756944,   "Bought (Third Party PD - 2) X 15",  -195.0),
  @(44864.9279996875,    "Bought (Third Party PD - 5) X 13",  -208.0),
  @(44864.928822766204,  "Bought (Third Party PD - 1) X 40",  -480.0)
)

$row = 16
foreach ($d in $newRows) {
    $ws1.Cells.Item($row, 1).Value = $d[0]
    $ws1.Cells.Item($row, 1).NumberFormat = "m/d/yyyy"
    $ws1.Cells.Item($row, 2).Value = $d[1]
    $ws1.Cells.Item($row, 3).Value = $d[2]
    $row = $row + 1
}

# Selection / active-sheet state recorded at save time.
$ws1.Range("A13:XFD13").Select()
$ws2.Activate()
